$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B12").Select()
